# Fix for shift swap in "intercambio de turnos" (shift exchange) table.
# The K:N block (rows 3-57) on sheet "Trabajador" is shifted up by one row:
# new row N (3..56) gets the old values of row N+1, and the last row (57)
# is cleared, since the trailing blank row (58) had no values either.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Trabajador")

$firstRow = 3
$lastRow = 57

# Capture original K:N values (columns 11..14) before overwriting anything.
$src = $ws.Range("K$($firstRow + 1):N$($lastRow + 1)").Value()

# Write them one row higher.
$destRowCount = $lastRow - $firstRow + 1
for ($i = 1; $i -le $destRowCount; $i++) {
    $destRow = $firstRow + $i - 1
    $ws.Cells.Item($destRow, 11).Value = $src[$i, 1]
    $ws.Cells.Item($destRow, 12).Value = $src[$i, 2]
    $ws.Cells.Item($destRow, 13).Value = $src[$i, 3]
    $ws.Cells.Item($destRow, 14).Value = $src[$i, 4]
}

# The former last row's values have now been duplicated into row 56; clear
# the vacated last row (57) so it matches the trailing blank rows below it.
$ws.Range("K$($lastRow):N$($lastRow)").ClearContents()
